$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 (I0) and J1 (IF) - reuse the header style from H1 (s="1")
# by copying H1's formatting, then set the text values.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF), rows 2-41
$data = @(
    @(4, 5),
    @(5, 5),
    @(6, 6),
    @(10, 10),
    @(7, 7),
    @(7, 7),
    @(7, 7),
    @(7, 7),
    @(10, 10),
    @(7, 7),
    @(7, 7),
    @(5, 5),
    @(7, 7),
    @(8, 8),
    @(10, 10),
    @(10, 10),
    @(8, 8),
    @(9, 9),
    @(7, 7),
    @(6, 7),
    @(7, 7),
    @(7, 7),
    @(6, 6),
    @(7, 7),
    @(8, 9),
    @(8, 8),
    @(6, 6),
    @(8, 8),
    @(8, 8),
    @(9, 9),
    @(1, 3),
    @(1, 2),
    @(5, 6),
    @(6, 6),
    @(6, 6),
    @(9, 9),
    @(7, 7),
    @(8, 8),
    @(7, 7),
    @(5, 5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}

$wb.Save()
